$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.000106692314148
$ws.Range("B1").Value = 0.916583776473999
$ws.Range("C1").Value = 0.7083269357681274
$ws.Range("D1").Value = 0.6849868297576904
$ws.Range("E1").Value = 0.7384288907051086
